# Update the "Number of Genotypes" column (column 6) of the outplant
# table: every data row's genotype count becomes "2", and the Totals
# row's genotype sum becomes "34" (17 rows x 2).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$lastRow = $t.Rows.Count   # 19: 1 header + 17 data rows + 1 totals row
$lastCol = 6                # "Number of Genotypes" column

# Data rows: row 2 .. (lastRow - 1)
for ($r = 2; $r -lt $lastRow; $r++) {
    $cell = $t.Cell($r, $lastCol)
    $cell.Range.Text = "2"
}

# Totals row: sum of all genotype counts (17 * 2 = 34)
$totalsCell = $t.Cell($lastRow, $lastCol)
$totalsCell.Range.Text = "34"
